# debug : [InvKinematics] : probrem = can't stack ik_data_st_ -> I generated and
# registered InvKinematics many times.
#
# Renumbers the existing IK rows (4-6), splits the old "elbow IK" row 7 into a
# separate wrist-rotation row, and fills in the two previously-blank rows
# (8-9) with the new wrist-IK (pos) entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 : RATA01 / RA01 -------------------------------------------------
$ws.Range("A4").Value = 2301
$ws.Range("B4").Value = "RATA01"
$ws.Range("C4").Value = 300
$ws.Range("D4").Value = "RA01"
$ws.Range("E4").Value = 1300
$ws.Range("F4").Value = "RAEE01"
$ws.Range("G4").Value = 2300
$ws.Range("H4").Value = "RATA01"
$ws.Range("I4").Value = "11: diry_look_pos"
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = "on"
$ws.Range("L4").Value = "off"
$ws.Range("M4").Value = "on"
$ws.Range("N4").Value = "Right wrist Rot IK"

# --- Row 5 : RATA02 / RA02 -------------------------------------------------
$ws.Range("A5").Value = 2302
$ws.Range("B5").Value = "RATA02"
$ws.Range("C5").Value = 301
$ws.Range("D5").Value = "RA02"
$ws.Range("E5").Value = 1300
$ws.Range("F5").Value = "RAEE01"
$ws.Range("G5").Value = 2300
$ws.Range("H5").Value = "RATA01"
$ws.Range("I5").Value = "11: diry_look_pos"
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = "on"
$ws.Range("L5").Value = "off"
$ws.Range("M5").Value = "on"
$ws.Range("N5").Value = "Right wrist Rot IK"

# --- Row 6 : RATA03 / RA03 -------------------------------------------------
$ws.Range("A6").Value = 2303
$ws.Range("B6").Value = "RATA03"
$ws.Range("C6").Value = 302
$ws.Range("D6").Value = "RA03"
$ws.Range("E6").Value = 1300
$ws.Range("F6").Value = "RAEE01"
$ws.Range("G6").Value = 2300
$ws.Range("H6").Value = "RATA01"
$ws.Range("I6").Value = "11: diry_look_pos"
$ws.Range("J6").Value = 0.2
$ws.Range("K6").Value = "on"
$ws.Range("L6").Value = "off"
$ws.Range("M6").Value = "on"
$ws.Range("N6").Value = "Right wrist Rot IK"

# --- Row 7 : RATA04 / RA01 (right elbow IK) --------------------------------
$ws.Range("A7").Value = 2304
$ws.Range("B7").Value = "RATA04"
$ws.Range("C7").Value = 300
$ws.Range("D7").Value = "RA01"
$ws.Range("E7").Value = 1301
$ws.Range("F7").Value = "RAEE02"
$ws.Range("G7").Value = 2301
$ws.Range("H7").Value = "RATA02"
$ws.Range("I7").Value = "0: pos_to_pos"
$ws.Range("J7").Value = 0.2
$ws.Range("K7").Value = "on"
$ws.Range("L7").Value = "off"
$ws.Range("M7").Value = "on"
$ws.Range("N7").Value = "Right elbow IK"

# --- Row 8 (was blank) : RATA05 / RA02 (right wrist Pos IK) ----------------
$ws.Range("A8").Value = 2305
$ws.Range("B8").Value = "RATA05"
$ws.Range("C8").Value = 301
$ws.Range("D8").Value = "RA02"
$ws.Range("E8").Value = 1302
$ws.Range("F8").Value = "RAEE03"
$ws.Range("G8").Value = 2302
$ws.Range("H8").Value = "RATA03"
$ws.Range("I8").Value = "0: pos_to_pos"
$ws.Range("J8").Value = 0.8
$ws.Range("K8").Value = "on"
$ws.Range("L8").Value = "off"
$ws.Range("M8").Value = "on"
$ws.Range("N8").Value = "Right wrist Pos IK"

# --- Row 9 (was blank) : RATA06 / RA01 (right wrist Pos IK) ----------------
# Row 9 column A previously used style 6 (plain) ; it now needs to match the
# bordered style used by the rest of this IK block (same as A8), so copy the
# formats over before writing the values.
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A9").Value = 2306
$ws.Range("B9").Value = "RATA06"
$ws.Range("C9").Value = 300
$ws.Range("D9").Value = "RA01"
$ws.Range("E9").Value = 1302
$ws.Range("F9").Value = "RAEE03"
$ws.Range("G9").Value = 2302
$ws.Range("H9").Value = "RATA03"
$ws.Range("I9").Value = "0: pos_to_pos"
$ws.Range("J9").Value = 0.1
$ws.Range("K9").Value = "on"
$ws.Range("L9").Value = "on"
$ws.Range("M9").Value = "on"
$ws.Range("N9").Value = "Right wrist Pos IK"

# Selection ends up on J8, matching the author's final cursor position.
$ws.Range("J8").Select()
